$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$arr = New-Object 'object[,]' 9,20
$arr[0,0] = "ECs"
$arr[0,1] = "Vtn"
$arr[0,2] = "Tnfrsf11b"
$arr[0,3] = "ECs"
$arr[0,4] = 3
$arr[0,5] = 1
$arr[0,6] = 5.835941000000001
$arr[0,7] = 17.507823
$arr[0,8] = 0.03643643319117328
$arr[0,9] = 0.03643643319117327
$arr[0,10] = 1
$arr[0,11] = 0.3333333333333333
$arr[0,12] = 0.07580833333333332
$arr[0,13] = 0.227425
$arr[0,14] = 0.028190957994264
$arr[0,15] = 0.02819095799426401
$arr[0,16] = 0.4424129606416667
$arr[0,17] = 3.981716645775
$arr[0,18] = 0.001027177957553173
$arr[0,19] = 0.001027177957553172
$arr[1,0] = "ECs"
$arr[1,1] = "Vtn"
$arr[1,2] = "Tnfrsf11b"
$arr[1,3] = "FAPs"
$arr[1,4] = 3
$arr[1,5] = 1
$arr[1,6] = 5.835941000000001
$arr[1,7] = 17.507823
$arr[1,8] = 0.03643643319117328
$arr[1,9] = 0.03643643319117327
$arr[1,10] = 3
$arr[1,11] = 1
$arr[1,12] = 2.095195666666667
$arr[1,13] = 6.285587
$arr[1,14] = 0.7791435378093522
$arr[1,15] = 0.7791435378093522
$arr[1,16] = 12.22743829412233
$arr[1,17] = 110.046944647101
$arr[1,18] = 0.02838921146172485
$arr[1,19] = 0.02838921146172484
$arr[2,0] = "ECs"
$arr[2,1] = "Vtn"
$arr[2,2] = "Tnfrsf11b"
$arr[2,3] = "MuSCs"
$arr[2,4] = 3
$arr[2,5] = 1
$arr[2,6] = 5.835941000000001
$arr[2,7] = 17.507823
$arr[2,8] = 0.03643643319117328
$arr[2,9] = 0.03643643319117327
$arr[2,10] = 3
$arr[2,11] = 1
$arr[2,12] = 0.5180969999999999
$arr[2,13] = 1.554291
$arr[2,14] = 0.1926655041963838
$arr[2,15] = 0.1926655041963838
$arr[2,16] = 3.023583524277
$arr[2,17] = 27.212251718493
$arr[2,18] = 0.007020043771895254
$arr[2,19] = 0.007020043771895253
$arr[3,0] = "FAPs"
$arr[3,1] = "Vtn"
$arr[3,2] = "Tnfrsf11b"
$arr[3,3] = "ECs"
$arr[3,4] = 3
$arr[3,5] = 1
$arr[3,6] = 17.50798033333334
$arr[3,7] = 52.52394100000001
$arr[3,8] = 0.1093102818770573
$arr[3,9] = 0.1093102818770573
$arr[3,10] = 1
$arr[3,11] = 0.3333333333333333
$arr[3,12] = 0.07580833333333332
$arr[3,13] = 0.227425
$arr[3,14] = 0.028190957994264
$arr[3,15] = 0.02819095799426401
$arr[3,16] = 1.327250809102778
$arr[3,17] = 11.945257281925
$arr[3,18] = 0.00308156156473728
$arr[3,19] = 0.00308156156473728
$arr[4,0] = "FAPs"
$arr[4,1] = "Vtn"
$arr[4,2] = "Tnfrsf11b"
$arr[4,3] = "FAPs"
$arr[4,4] = 3
$arr[4,5] = 1
$arr[4,6] = 17.50798033333334
$arr[4,7] = 52.52394100000001
$arr[4,8] = 0.1093102818770573
$arr[4,9] = 0.1093102818770573
$arr[4,10] = 3
$arr[4,11] = 1
$arr[4,12] = 2.095195666666667
$arr[4,13] = 6.285587
$arr[4,14] = 0.7791435378093522
$arr[4,15] = 0.7791435378093522
$arr[4,16] = 36.68264452648523
$arr[4,17] = 330.143800738367
$arr[4,18] = 0.08516839974062794
$arr[4,19] = 0.08516839974062793
$arr[5,0] = "FAPs"
$arr[5,1] = "Vtn"
$arr[5,2] = "Tnfrsf11b"
$arr[5,3] = "MuSCs"
$arr[5,4] = 3
$arr[5,5] = 1
$arr[5,6] = 17.50798033333334
$arr[5,7] = 52.52394100000001
$arr[5,8] = 0.1093102818770573
$arr[5,9] = 0.1093102818770573
$arr[5,10] = 3
$arr[5,11] = 1
$arr[5,12] = 0.5180969999999999
$arr[5,13] = 1.554291
$arr[5,14] = 0.1926655041963838
$arr[5,15] = 0.1926655041963838
$arr[5,16] = 9.070832086758999
$arr[5,17] = 81.63748878083101
$arr[5,18] = 0.02106032057169208
$arr[5,19] = 0.02106032057169208
$arr[6,0] = "MuSCs"
$arr[6,1] = "Vtn"
$arr[6,2] = "Tnfrsf11b"
$arr[6,3] = "ECs"
$arr[6,4] = 3
$arr[6,5] = 1
$arr[6,6] = 136.8238143333333
$arr[6,7] = 410.471443
$arr[6,8] = 0.8542532849317694
$arr[6,9] = 0.8542532849317694
$arr[6,10] = 1
$arr[6,11] = 0.3333333333333333
$arr[6,12] = 0.07580833333333332
$arr[6,13] = 0.227425
$arr[6,14] = 0.028190957994264
$arr[6,15] = 0.02819095799426401
$arr[6,16] = 10.37238532491944
$arr[6,17] = 93.351467924275
$arr[6,18] = 0.02408221847197355
$arr[6,19] = 0.02408221847197355
$arr[7,0] = "MuSCs"
$arr[7,1] = "Vtn"
$arr[7,2] = "Tnfrsf11b"
$arr[7,3] = "FAPs"
$arr[7,4] = 3
$arr[7,5] = 1
$arr[7,6] = 136.8238143333333
$arr[7,7] = 410.471443
$arr[7,8] = 0.8542532849317694
$arr[7,9] = 0.8542532849317694
$arr[7,10] = 3
$arr[7,11] = 1
$arr[7,12] = 2.095195666666667
$arr[7,13] = 6.285587
$arr[7,14] = 0.7791435378093522
$arr[7,15] = 0.7791435378093522
$arr[7,16] = 286.6726628880045
$arr[7,17] = 2580.053965992041
$arr[7,18] = 0.6655859266069994
$arr[7,19] = 0.6655859266069994
$arr[8,0] = "MuSCs"
$arr[8,1] = "Vtn"
$arr[8,2] = "Tnfrsf11b"
$arr[8,3] = "MuSCs"
$arr[8,4] = 3
$arr[8,5] = 1
$arr[8,6] = 136.8238143333333
$arr[8,7] = 410.471443
$arr[8,8] = 0.8542532849317694
$arr[8,9] = 0.8542532849317694
$arr[8,10] = 3
$arr[8,11] = 1
$arr[8,12] = 0.5180969999999999
$arr[8,13] = 1.554291
$arr[8,14] = 0.1926655041963838
$arr[8,15] = 0.1926655041963838
$arr[8,16] = 70.88800773465699
$arr[8,17] = 637.9920696119129
$arr[8,18] = 0.1645851398527965
$arr[8,19] = 0.1645851398527965

$rng = $ws.Range("A2:T10")
$rng.Value = $arr
